{"js": "// Rewrite the bullet list under \"KEY ACHIEVEMENTS AND IMPACT\" (the \"Impact\"\n// sub-heading) to the new, impact-focused accomplishment statements, and\n// drop the two extra bullets that are no longer wanted.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the \"KEY ACHIEVEMENTS AND IMPACT\" Heading 2 paragraph.\nlet sectionIdx = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.style === \"Heading 2\" && p.text.trim() === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    sectionIdx = i;\n    break;\n  }\n}\nif (sectionIdx === -1) {\n  throw new Error('Could not find \"KEY ACHIEVEMENTS AND IMPACT\" heading');\n}\n\n// Collect the bullet paragraphs that belong to that section: everything\n// after the heading (and its \"Impact\" sub-heading) up to the next\n// Heading 1/2 paragraph.\nconst bulletParas = [];\nfor (let i = sectionIdx + 1; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.style === \"Heading 1\" || p.style === \"Heading 2\") {\n    break;\n  }\n  if (p.style === \"Heading 3\") {\n    continue; // the \"Impact\" sub-heading\n  }\n  bulletParas.push(p);\n}\n\nconst newBullets = [\n  \"\u2022 Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard\",\n  \"\u2022 Reduced polling margins from \u00b14.2% to \u00b12.1%\",\n  \"\u2022 Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis\",\n  \"\u2022 Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\"\n];\n\n// Replace the text of the first four bullet paragraphs in place.\nfor (let i = 0; i < newBullets.length; i++) {\n  bulletParas[i].insertText(newBullets[i], Word.InsertLocation.replace);\n}\n\n// Delete the remaining (now unwanted) bullet paragraphs, in reverse order\n// so earlier deletions don't shift the indices of the ones still to remove.\nfor (let i = bulletParas.length - 1; i >= newBullets.length; i--) {\n  bulletParas[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Rewrite the bullet list under \"KEY ACHIEVEMENTS AND IMPACT\" (the \"Impact\"\n# sub-heading) to the new, impact-focused accomplishment statements, and\n# drop the two extra bullets that are no longer wanted.\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n# Locate the \"KEY ACHIEVEMENTS AND IMPACT\" Heading 2 paragraph (1-based COM index).\n$sectionIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Style.NameLocal -eq \"Heading 2\" -and $p.Range.Text.Trim() -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $sectionIdx = $i\n        break\n    }\n}\nif ($sectionIdx -eq -1) {\n    throw \"Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading\"\n}\n\n# Collect the indices of the bullet paragraphs belonging to that section:\n# everything after the heading (and its \"Impact\" sub-heading) up to the\n# next Heading 1/2 paragraph.\n$bulletIdx = @()\nfor ($i = $sectionIdx + 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    $styleName = $p.Style.NameLocal\n    if ($styleName -eq \"Heading 1\" -or $styleName -eq \"Heading 2\") {\n        break\n    }\n    if ($styleName -eq \"Heading 3\") {\n        continue\n    }\n    $bulletIdx += $i\n}\n\n$newBullets = @(\n    \"\u2022 Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard\",\n    \"\u2022 Reduced polling margins from \u00b14.2% to \u00b12.1%\",\n    \"\u2022 Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis\",\n    \"\u2022 Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\"\n)\n\n# Replace the text of the first four bullet paragraphs in place.\nfor ($k = 0; $k -lt $newBullets.Length; $k++) {\n    $paras.Item($bulletIdx[$k]).Range.Text = $newBullets[$k]\n}\n\n# Delete the remaining (now unwanted) bullet paragraphs, in reverse order\n# so earlier deletions don't shift the indices of the ones still to remove.\nfor ($k = $bulletIdx.Length - 1; $k -ge $newBullets.Length; $k--) {\n    $paras.Item($bulletIdx[$k]).Range.Delete()\n}\n"}
